$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Add new rows 5 & 6 to Sheet2 (20201001-Actin then 20200925-Actin)
$ws2.Range("B5").Value = "F:\PhD, PMMH, ESPCI\Processing\20201001-Actin\results\Classification manually 20201001-Actin.xlsx"
$ws2.Range("C5").Value = "E:\Dropbox\Research\All Plottings\20201001-Actin\Figures"

$ws2.Range("B6").Value = "F:\PhD, PMMH, ESPCI\Processing\20200925-Actin\results\Classification manually 20200925-Actin.xlsx"
$ws2.Range("C6").Value = "E:\Dropbox\Research\All Plottings\20200925-Actin\Figures"

# Update the selections left behind on each sheet
$ws1.Activate()
$ws1.Range("D19").Select()

$ws2.Activate()
$ws2.Range("B8").Select()
